$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "ChangePassword"

$ws.Range("A1").Value = "NewPassword"
$ws.Range("B1").Value = "ConfirmPassword"
$ws.Range("A2").Value = "admin@123456"
$ws.Range("B2").Value = "admin@123456"

$ws.Range("A1:B2").Borders.LineStyle = 1
$ws.Range("A1:B2").Borders.Weight = 2

$ws.Range("A1:B1").Interior.ThemeColor = 4
$ws.Range("A1:B1").Interior.TintAndShade = 0.59999389629810485

